# Update "想去人数" (column F) values for rows 2-14 on the
# "展览" and "全部类型" worksheets, matching the latest generated output.

$wb = $excel.ActiveWorkbook

$newValues = @{
    2  = 1617
    3  = 214
    4  = 203
    5  = 6138
    6  = 357
    7  = 242
    8  = 50
    9  = 13
    10 = 8843
    11 = 2359
    12 = 257
    13 = 5594
    14 = 10317
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Cells.Item($row, 6).Value = $newValues[$row]
    }
}
